$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").Value = "Mail Analytics Workshop"
$ws.Range("E5").Select()
